$wb = $excel.ActiveWorkbook

# ---- Update sigma_010 (sheet index 2) ----
$ws2 = $wb.Worksheets.Item("sigma_010")
$ws2.Cells.Item(2, 2).Value = 27.80235843975257
$ws2.Cells.Item(2, 3).Value = 30.08023247722995
$ws2.Cells.Item(3, 2).Value = 27.84526899161106
$ws2.Cells.Item(3, 3).Value = 30.09481007481382
$ws2.Cells.Item(4, 2).Value = 27.81797992354754
$ws2.Cells.Item(4, 3).Value = 30.09434967994184
$ws2.Cells.Item(5, 2).Value = 27.82290678197702
$ws2.Cells.Item(5, 3).Value = 30.08852874083704
$ws2.Cells.Item(6, 2).Value = 27.81343603938643
$ws2.Cells.Item(6, 3).Value = 30.09827504513417
$ws2.Cells.Item(7, 2).Value = 27.82527888409643
$ws2.Cells.Item(7, 3).Value = 30.09968142856989
$ws2.Cells.Item(8, 2).Value = 27.81139993429294
$ws2.Cells.Item(8, 3).Value = 30.08512053339445
$ws2.Cells.Item(9, 2).Value = 27.80757204591135
$ws2.Cells.Item(9, 3).Value = 30.08638560087784
$ws2.Cells.Item(10, 2).Value = 27.85891761061971
$ws2.Cells.Item(10, 3).Value = 30.12175982841377
$ws2.Cells.Item(11, 2).Value = 27.78739978112536
$ws2.Cells.Item(11, 3).Value = 30.07130906341056
$ws2.Cells.Item(12, 2).Value = 27.81925184323204
$ws2.Cells.Item(12, 3).Value = 30.09204524726233

# ---- Update sigma_025 (sheet index 3) ----
$ws3 = $wb.Worksheets.Item("sigma_025")
$ws3.Cells.Item(2, 2).Value = 19.76259059450875
$ws3.Cells.Item(2, 3).Value = 26.46205091117629
$ws3.Cells.Item(3, 2).Value = 19.72653186943145
$ws3.Cells.Item(3, 3).Value = 26.49740113483317
$ws3.Cells.Item(4, 2).Value = 19.73922447521184
$ws3.Cells.Item(4, 3).Value = 26.49283122699845
$ws3.Cells.Item(5, 2).Value = 19.734632245899
$ws3.Cells.Item(5, 3).Value = 26.4747809489479
$ws3.Cells.Item(6, 2).Value = 19.73500410834658
$ws3.Cells.Item(6, 3).Value = 26.48748043226867
$ws3.Cells.Item(7, 2).Value = 19.75227152963232
$ws3.Cells.Item(7, 3).Value = 26.48106842000163
$ws3.Cells.Item(8, 2).Value = 19.74549604676852
$ws3.Cells.Item(8, 3).Value = 26.48745651417626
$ws3.Cells.Item(9, 2).Value = 19.73165605984333
$ws3.Cells.Item(9, 3).Value = 26.48236437996426
$ws3.Cells.Item(10, 2).Value = 19.73966720961584
$ws3.Cells.Item(10, 3).Value = 26.52882112959994
$ws3.Cells.Item(11, 2).Value = 19.74509029740211
$ws3.Cells.Item(11, 3).Value = 26.50559768232268
$ws3.Cells.Item(12, 2).Value = 19.74121644366597
$ws3.Cells.Item(12, 3).Value = 26.48998527802893

# ---- Add new sheet sigma_050 after sigma_025 ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "sigma_050"

$ws4.Cells.Item(1, 1).Value = "Rows"
$ws4.Cells.Item(1, 2).Value = "Noisy"
$ws4.Cells.Item(1, 3).Value = "NLM-LBP"

$ws4.Cells.Item(2, 1).Value = 0
$ws4.Cells.Item(2, 2).Value = 14.80041056161648
$ws4.Cells.Item(2, 3).Value = 21.33534017853955
$ws4.Cells.Item(3, 1).Value = 1
$ws4.Cells.Item(3, 2).Value = 14.79375654584836
$ws4.Cells.Item(3, 3).Value = 21.33436801443766
$ws4.Cells.Item(4, 1).Value = 2
$ws4.Cells.Item(4, 2).Value = 14.78789335695933
$ws4.Cells.Item(4, 3).Value = 21.29837443538272
$ws4.Cells.Item(5, 1).Value = 3
$ws4.Cells.Item(5, 2).Value = 14.79715515651936
$ws4.Cells.Item(5, 3).Value = 21.36288016426221
$ws4.Cells.Item(6, 1).Value = 4
$ws4.Cells.Item(6, 2).Value = 14.79541971734821
$ws4.Cells.Item(6, 3).Value = 21.281727562337
$ws4.Cells.Item(7, 1).Value = 5
$ws4.Cells.Item(7, 2).Value = 14.78248594624705
$ws4.Cells.Item(7, 3).Value = 21.30261495800624
$ws4.Cells.Item(8, 1).Value = 6
$ws4.Cells.Item(8, 2).Value = 14.7724845914075
$ws4.Cells.Item(8, 3).Value = 21.30132892711007
$ws4.Cells.Item(9, 1).Value = 7
$ws4.Cells.Item(9, 2).Value = 14.77873068294177
$ws4.Cells.Item(9, 3).Value = 21.29598751135684
$ws4.Cells.Item(10, 1).Value = 8
$ws4.Cells.Item(10, 2).Value = 14.79168311475186
$ws4.Cells.Item(10, 3).Value = 21.32543175175758
$ws4.Cells.Item(11, 1).Value = 9
$ws4.Cells.Item(11, 2).Value = 14.78027577375493
$ws4.Cells.Item(11, 3).Value = 21.29265490057226
$ws4.Cells.Item(12, 1).Value = "Média"
$ws4.Cells.Item(12, 2).Value = 14.78802954473948
$ws4.Cells.Item(12, 3).Value = 21.31307084037621
